$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1424953361892989
$ws.Range("D2").Value = 0.0004299426552520913
$ws.Range("E2").Value = 0.07991781540580511
$ws.Range("F2").Value = 0.5932467633878105
$ws.Range("G2").Value = 0.002345679738735331
$ws.Range("M2").Value = 0.6979156703716711
$ws.Range("N2").Value = 1.437313047148649
$ws.Range("O2").Value = 1.897006763755314
$ws.Range("B3").Value = 0.1329157298043668
$ws.Range("D3").Value = 0.0003964902536868209
$ws.Range("E3").Value = 0.08287708702156316
$ws.Range("F3").Value = 0.5571254389356568
$ws.Range("G3").Value = 0.002349592272756002
$ws.Range("M3").Value = 0.6114803225277541
$ws.Range("N3").Value = 1.386698451232718
$ws.Range("O3").Value = 1.784941032402287
$ws.Range("B4").Value = 0.1271001053738843
$ws.Range("D4").Value = 0.0003761996499376252
$ws.Range("E4").Value = 0.08484063289508459
$ws.Range("F4").Value = 0.5353679329008258
$ws.Range("G4").Value = 0.00235212330709289
$ws.Range("M4").Value = 0.5583235224461163
$ws.Range("N4").Value = 1.356146371068775
$ws.Range("O4").Value = 1.717523497416835
$ws.Range("B5").Value = 0.1247470728135198
$ws.Range("D5").Value = 0.0003679924427686387
$ws.Range("E5").Value = 0.08567735536146248
$ws.Range("F5").Value = 0.5266066972169909
$ws.Range("G5").Value = 0.002353187199389955
$ws.Range("M5").Value = 0.5366408634389899
$ws.Range("N5").Value = 1.343830187433582
$ws.Range("O5").Value = 1.690397430685408
$ws.Range("B6").Value = 0.124357379393004
$ws.Range("D6").Value = 0.0003666333134066946
$ws.Range("E6").Value = 0.08581849232594951
$ws.Range("F6").Value = 0.5251582274959219
$ws.Range("G6").Value = 0.00235336582243832
$ws.Range("M6").Value = 0.5330392336333318
$ws.Range("N6").Value = 1.341793246235113
$ws.Range("O6").Value = 1.685914054317834
$ws.Range("B7").Value = 0.1270683029193833
$ws.Range("D7").Value = 0.0003760887179913297
$ws.Range("E7").Value = 0.08485176957115215
$ws.Range("F7").Value = 0.5352493511896341
$ws.Range("G7").Value = 0.002352137523606063
$ws.Range("M7").Value = 0.5580311860887832
$ws.Range("N7").Value = 1.355979725467307
$ws.Range("O7").Value = 1.717156263897948
$ws.Range("B8").Value = 0.139178652560517
$ws.Range("D8").Value = 0.0004183556874206573
$ws.Range("E8").Value = 0.08090758485342597
$ws.Range("F8").Value = 0.5807042463160599
$ws.Range("G8").Value = 0.002347002128521032
$ws.Range("M8").Value = 0.668130752854978
$ws.Range("N8").Value = 1.419753184325828
$ws.Range("O8").Value = 1.858076060527935
$ws.Range("B9").Value = 0.1634450174328919
$ws.Range("D9").Value = 0.0005032886288409699
$ws.Range("E9").Value = 0.07434860355879103
$ws.Range("F9").Value = 0.6732227270668432
$ws.Range("G9").Value = 0.002337948122705017
$ws.Range("M9").Value = 0.8833456807207654
$ws.Range("N9").Value = 1.548906929911681
$ws.Range("O9").Value = 2.145597535330467
$ws.Range("B10").Value = 0.1815804850218825
$ws.Range("D10").Value = 0.0005670481755615242
$ws.Range("E10").Value = 0.07026389571659841
$ws.Range("F10").Value = 0.7433205311375986
$ws.Range("G10").Value = 0.002331909010651221
$ws.Range("M10").Value = 1.041040522520674
$ws.Range("N10").Value = 1.646202818567332
$ws.Range("O10").Value = 2.36386996329054
$ws.Range("B11").Value = 0.1898956811295278
$ws.Range("D11").Value = 0.0005963762129219674
$ws.Range("E11").Value = 0.06856886495660497
$ws.Range("F11").Value = 0.775685224190596
$ws.Range("G11").Value = 0.002329293288163067
$ws.Range("M11").Value = 1.112688615605535
$ws.Range("N11").Value = 1.6909699208951
$ws.Range("O11").Value = 2.464742455019348
$ws.Range("B12").Value = 0.1930536271213015
$ws.Range("D12").Value = 0.0006075307549693321
$ws.Range("E12").Value = 0.06795077603490007
$ws.Range("F12").Value = 0.7880104474648419
$ws.Range("G12").Value = 0.002328321580814079
$ws.Range("M12").Value = 1.139806921705983
$ws.Range("N12").Value = 1.707993225147618
$ws.Range("O12").Value = 2.503170763403546
$ws.Range("B13").Value = 0.192373103071688
$ws.Range("D13").Value = 0.0006051262258424117
$ws.Range("E13").Value = 0.06808282960563439
$ws.Range("F13").Value = 0.7853528892544688
$ws.Range("G13").Value = 0.002328530020276976
$ws.Range("M13").Value = 1.133967109040384
$ws.Range("N13").Value = 1.704323825273605
$ws.Range("O13").Value = 2.49488425870868
$ws.Range("B14").Value = 0.1901553049249003
$ws.Range("D14").Value = 0.0005972929162041396
$ws.Range("E14").Value = 0.06851753590111898
$ws.Range("F14").Value = 0.7766978316513047
$ws.Range("G14").Value = 0.002329212968710288
$ws.Range("M14").Value = 1.114919924375073
$ws.Range("N14").Value = 1.692369026707382
$ws.Range("O14").Value = 2.467899346638205
$ws.Range("B15").Value = 0.1887980259508311
$ws.Range("D15").Value = 0.0005925011900451693
$ws.Range("E15").Value = 0.06878691293332473
$ws.Range("F15").Value = 0.7714054242384378
$ws.Range("G15").Value = 0.002329633741001528
$ws.Range("M15").Value = 1.103251227676083
$ws.Range("N15").Value = 1.6850555545806
$ws.Range("O15").Value = 2.451400361737797
$ws.Range("B16").Value = 0.1810383617286391
$ws.Range("D16").Value = 0.0005651382078308842
$ws.Range("E16").Value = 0.07037798019597474
$ws.Range("F16").Value = 0.7412150899513392
$ws.Range("G16").Value = 0.002332082591947502
$ws.Range("M16").Value = 1.036356326608924
$ws.Range("N16").Value = 1.643287217421204
$ws.Range("O16").Value = 2.357309742578934
$ws.Range("B17").Value = 0.1762946184794174
$ws.Range("D17").Value = 0.0005484362326200198
$ws.Range("E17").Value = 0.07139605912582248
$ws.Range("F17").Value = 0.722817030495122
$ws.Range("G17").Value = 0.002333618491805741
$ws.Range("M17").Value = 0.9952954763112416
$ws.Range("N17").Value = 1.617792125056582
$ws.Range("O17").Value = 2.299994803115737
$ws.Range("B18").Value = 0.1735723008559233
$ws.Range("D18").Value = 0.0005388600111704278
$ws.Range("E18").Value = 0.07199696554792467
$ws.Range("F18").Value = 0.7122797857831813
$ws.Range("G18").Value = 0.002334514284262261
$ws.Range("M18").Value = 0.9716701059242894
$ws.Range("N18").Value = 1.603175828286851
$ws.Range("O18").Value = 2.267177234223823
$ws.Range("B19").Value = 0.1726516354098635
$ws.Range("D19").Value = 0.000535622810675207
$ws.Range("E19").Value = 0.07220304594288862
$ws.Range("F19").Value = 0.7087197301220698
$ws.Range("G19").Value = 0.002334819713929099
$ws.Range("M19").Value = 0.9636695622874072
$ws.Range("N19").Value = 1.598235270037293
$ws.Range("O19").Value = 2.256091178172198
$ws.Range("B20").Value = 0.1767989623493804
$ws.Range("D20").Value = 0.0005502110340138699
$ws.Range("E20").Value = 0.07128609349157955
$ws.Range("F20").Value = 0.7247708892399061
$ws.Range("G20").Value = 0.002333453711724887
$ws.Range("M20").Value = 0.9996673327339636
$ws.Range("N20").Value = 1.620501187190825
$ws.Range("O20").Value = 2.306080688805764
$ws.Range("B21").Value = 0.190806479151945
$ws.Range("D21").Value = 0.0005995924112927042
$ws.Range("E21").Value = 0.06838920407091642
$ws.Range("F21").Value = 0.7792381411566538
$ws.Range("G21").Value = 0.00232901186062402
$ws.Range("M21").Value = 1.1205149091723
$ws.Range("N21").Value = 1.695878529578522
$ws.Range("O21").Value = 2.47581920107217
$ws.Range("B22").Value = 0.2000144978829184
$ws.Range("D22").Value = 0.0006321506626214557
$ws.Range("E22").Value = 0.06663469428964142
$ws.Range("F22").Value = 0.8152406480043481
$ws.Range("G22").Value = 0.002326218447944428
$ws.Range("M22").Value = 1.199418157179153
$ws.Range("N22").Value = 1.745554809723473
$ws.Range("O22").Value = 2.588095645181738
$ws.Range("B23").Value = 0.1950952028399655
$ws.Range("D23").Value = 0.0006147469224728752
$ws.Range("E23").Value = 0.06755830313418798
$ws.Range("F23").Value = 0.7959880930480807
$ws.Range("G23").Value = 0.002327699349673376
$ws.Range("M23").Value = 1.157313308303145
$ws.Range("N23").Value = 1.719004489371116
$ws.Range("O23").Value = 2.528047760697177
$ws.Range("B24").Value = 0.1765709331320835
$ws.Range("D24").Value = 0.000549408566016929
$ws.Range("E24").Value = 0.07133576038953926
$ws.Range("F24").Value = 0.7238874250313643
$ws.Range("G24").Value = 0.002333528169105112
$ws.Range("M24").Value = 0.9976908754013039
$ws.Range("N24").Value = 1.619276291812184
$ws.Range("O24").Value = 2.303328844109672
$ws.Range("B25").Value = 0.1568257795625811
$ws.Range("D25").Value = 0.0004800818385852779
$ws.Range("E25").Value = 0.07599532664092834
$ws.Range("F25").Value = 0.6478255176706256
$ws.Range("G25").Value = 0.002340289355582686
$ws.Range("M25").Value = 0.8251985613386239
$ws.Range("N25").Value = 1.513538423053916
$ws.Range("O25").Value = 2.503170763403546
